# Updates the cryptos list: refreshed price/volume figures plus a
# Quant <-> NEARProtocol row swap (rows 48-49), per commit:
# "Updated cryptos list on Thu May 18 02:06:04 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds plain numeric-looking text (e.g. "0.3752", "1.001").
# Force it to Text format before writing so Excel doesn't silently convert
# these into floating point numbers (which would drop things like trailing
# zeros / change precision). Formatting is cleared again afterwards so the
# cells end up with no explicit style, matching the original workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.353.46"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "1.822.50"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "314.13"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D7").Value = "0.4451"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("D8").Value = "0.3752"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "0.07471"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "0.8833"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").Value = "20.98"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "1.826.67"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "6.745"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "5.404"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "93.49"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "0.07122"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "0.000008755"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "15.12"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "27.358.66"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "5.399"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "2.052.47"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "1.960"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "151.29"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "2.319"
$ws.Range("E27").Value = "  +3.17%  "
$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "5.347"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "0.08898"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "0.7803"
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("D33").Value = "1.200"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "4.626"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").Value = "2.907"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "1.111"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "0.01991"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "0.05293"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "7.293"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "0.5302"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").Value = "2.857"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "0.1710"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "2.285"
$ws.Range("E44").Value = "  +14.94%  "
$ws.Range("D45").Value = "8.629"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "0.5073"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "10.52"
$ws.Range("E47").Value = "  -1.69%  "

# Row 48 becomes NEARProtocol, row 49 becomes Quant (swapped order).
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.691"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "105.06"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "0.06397"
$ws.Range("E51").Value = "  +0.59%  "

# Clear the temporary Text number format so cells return to the workbook's
# default (unstyled) formatting, matching the original authoring.
$priceRange.ClearFormats()
